# Add "competency" column to the lesson plan upload template.
# Inserts a new column before column F (which currently holds
# "learning_outcomes"), shifting the rest of the header row right by one,
# and fills it in with the header "competency" and the same "Number"
# placeholder value used by the other template columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F; everything from F onward shifts right.
$ws.Columns("F:F").Insert()

# Match the width of the column it was inserted in front of (18 chars,
# which Excel's ColumnWidth API reports as ~17.17 using this font's
# character-width metrics).
$ws.Columns("F:F").ColumnWidth = 17.17

# New header cell (old F1 = "learning_outcomes" is now G1).
$ws.Range("F1").Value = "competency"

# Row 2 placeholder value, matching the other template columns.
$ws.Range("F2").Value = "Number"

# Update the active selection to match the authored file.
$ws.Range("G6").Select()
